# Acc_fijos_N.xlsx edit script
# 1) Fix typo in the merged title D4: "Accesos a internte fijo" -> "Accesos a internet fijo"
# 2) Convert column C (Mes) values from numeric month (1-12) to Spanish month abbreviations
# 3) Adjust column C width to fit the new text content

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Replace numeric month values in column C with Spanish month abbreviations ---
$monthNames = @{
    1  = "Ene."
    2  = "Feb."
    3  = "Mar."
    4  = "Abr."
    5  = "May."
    6  = "Jun."
    7  = "Jul."
    8  = "Ago."
    9  = "Sep."
    10 = "Oct."
    11 = "Nov."
    12 = "Dic."
}

$monthOrder = @(8, 7, 6, 5, 4, 3, 2, 1, 12, 11, 10, 9)

foreach ($m in $monthOrder) {
    for ($r = 6; $r -le 80; $r++) {
        $cell = $ws.Cells.Item($r, 3)
        $cur = [int]$cell.Value2
        if ($cur -eq $m) {
            $cell.Value = $monthNames[$m]
        }
    }
}

# --- 2) Fix the typo in the merged header cell D4 ---
$ws.Range("D4").Value = "Accesos a internet fijo"

# --- 3) Resize column C so the new text fits (previously auto best-fit at width 4) ---
$ws.Columns("C").ColumnWidth = 5.7
